# "Working on Emp Education" — QC Test.xlsx updates
#
# 1. Model sheet (HR_EMP_DETAIL row): fill in the previously-blank
#    "Foreign Key [Fluent API]" cell with "No".
# 2. Controller sheet (EmpDetailController row): fill in the remaining
#    blank flags — "Yes" for Insert/Update/EQResult Messages/ModelState
#    Errors/DropDown, "No" for Delete.
# 3. Views sheet (EmpDetail row): fill in the whole blank row of flags.
# 4. The active sheet moved from "Controller" back to "Model", and the
#    selected cell on each sheet reflects where the author was working.

$wb = $excel.ActiveWorkbook

$modelWs      = $wb.Worksheets.Item("Model")
$controllerWs = $wb.Worksheets.Item("Controller")
$viewsWs      = $wb.Worksheets.Item("Views")

# --- Model sheet: HR_EMP_DETAIL row (row 8) ---
$modelWs.Range("D8").Value = "No"

# --- Controller sheet: EmpDetailController row (row 8) ---
$controllerWs.Range("G8").Value = "Yes"
$controllerWs.Range("H8").Value = "Yes"
$controllerWs.Range("I8").Value = "Yes"
$controllerWs.Range("J8").Value = "Yes"
$controllerWs.Range("K8").Value = "Yes"
$controllerWs.Range("L8").Value = "No"

# --- Views sheet: EmpDetail row (row 8) ---
$viewsWs.Range("B8").Value = "Yes"
$viewsWs.Range("C8").Value = "Yes"
$viewsWs.Range("D8").Value = "Yes"
$viewsWs.Range("E8").Value = "Yes"
$viewsWs.Range("F8").Value = "No"
$viewsWs.Range("G8").Value = "No"
$viewsWs.Range("H8").Value = "No"
$viewsWs.Range("I8").Value = "Yes"
$viewsWs.Range("J8").Value = "Yes"
$viewsWs.Range("K8").Value = "Yes"
$viewsWs.Range("L8").Value = "No"
$viewsWs.Range("M8").Value = "No"
$viewsWs.Range("N8").Value = "No"

# --- Update selections to match where the author ended up ---
$controllerWs.Range("L9").Select()
$viewsWs.Range("F12").Select()

# Re-select Model last and activate it so it becomes the active sheet/tab
# (matches removal of Controller's tabSelected + workbook activeTab moving
# back to the first sheet).
$modelWs.Activate()
$modelWs.Range("E10").Select()
